$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8781003333333334
$ws.Range("H2").Value = 2.634301
$ws.Range("I2").Value = 0.1010434633250494
$ws.Range("J2").Value = 0.1010434633250494
$ws.Range("M2").Value = 0.002913
$ws.Range("N2").Value = 0.008739
$ws.Range("Q2").Value = 0.002557906271
$ws.Range("R2").Value = 0.023021156439
$ws.Range("S2").Value = 0.1010434633250494
$ws.Range("T2").Value = 0.1010434633250494
$ws.Range("I3").Value = 0.01326751606355713
$ws.Range("J3").Value = 0.01326751606355713
$ws.Range("M3").Value = 0.002913
$ws.Range("N3").Value = 0.008739
$ws.Range("Q3").Value = 0.000335865987
$ws.Range("R3").Value = 0.003022793883
$ws.Range("S3").Value = 0.01326751606355713
$ws.Range("T3").Value = 0.01326751606355713
$ws.Range("G4").Value = 3.520787
$ws.Range("H4").Value = 10.562361
$ws.Range("I4").Value = 0.4051387963370292
$ws.Range("J4").Value = 0.4051387963370292
$ws.Range("M4").Value = 0.002913
$ws.Range("N4").Value = 0.008739
$ws.Range("Q4").Value = 0.010256052531
$ws.Range("R4").Value = 0.092304472779
$ws.Range("S4").Value = 0.4051387963370292
$ws.Range("T4").Value = 0.4051387963370292
$ws.Range("G5").Value = 0.0464
$ws.Range("H5").Value = 0.1392
$ws.Range("I5").Value = 0.005339272199663925
$ws.Range("J5").Value = 0.005339272199663925
$ws.Range("M5").Value = 0.002913
$ws.Range("N5").Value = 0.008739
$ws.Range("Q5").Value = 0.0001351632
$ws.Range("R5").Value = 0.0012164688
$ws.Range("S5").Value = 0.005339272199663925
$ws.Range("T5").Value = 0.005339272199663925
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.41099
$ws.Range("H6").Value = 1.23297
$ws.Range("I6").Value = 0.0472928336495663
$ws.Range("J6").Value = 0.0472928336495663
$ws.Range("M6").Value = 0.002913
$ws.Range("N6").Value = 0.008739
$ws.Range("Q6").Value = 0.00119721387
$ws.Range("R6").Value = 0.01077492483
$ws.Range("S6").Value = 0.0472928336495663
$ws.Range("T6").Value = 0.0472928336495663
$ws.Range("G7").Value = 3.718746666666667
$ws.Range("H7").Value = 11.15624
$ws.Range("I7").Value = 0.4279181184251341
$ws.Range("J7").Value = 0.4279181184251342
$ws.Range("M7").Value = 0.002913
$ws.Range("N7").Value = 0.008739
$ws.Range("Q7").Value = 0.01083270904
$ws.Range("R7").Value = 0.09749438136000001
$ws.Range("S7").Value = 0.4279181184251341
$ws.Range("T7").Value = 0.4279181184251342
